# Update NATMI LR-pair stats (Ebi3-Il6st) per Dr Hou's advice:
# ligand/receptor "expressing cells" counts go from 1 to 3, and the
# dependent average/total expression + specificity columns are
# recomputed accordingly for every data row (rows 2-11).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 3.143381999999999
$ws.Range("H2").Value = 9.430145999999999
$ws.Range("I2").Value = 0.3901511487518624
$ws.Range("J2").Value = 0.3901511487518625
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 86.40747066666665
$ws.Range("N2").Value = 259.222412
$ws.Range("O2").Value = 0.3380062309947018
$ws.Range("P2").Value = 0.3380062309947017
$ws.Range("Q2").Value = 271.6116879591279
$ws.Range("R2").Value = 2444.505191632151
$ws.Range("S2").Value = 0.1318735193078703
$ws.Range("T2").Value = 0.1318735193078702
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 3.143381999999999
$ws.Range("H3").Value = 9.430145999999999
$ws.Range("I3").Value = 0.3901511487518624
$ws.Range("J3").Value = 0.3901511487518625
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 125.002688
$ws.Range("N3").Value = 375.008064
$ws.Range("O3").Value = 0.4889818797969518
$ws.Range("P3").Value = 0.4889818797969517
$ws.Range("Q3").Value = 392.9311994108159
$ws.Range("R3").Value = 3536.380794697343
$ws.Range("S3").Value = 0.1907768421216259
$ws.Range("T3").Value = 0.1907768421216259
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 3.143381999999999
$ws.Range("H4").Value = 9.430145999999999
$ws.Range("I4").Value = 0.3901511487518624
$ws.Range("J4").Value = 0.3901511487518625
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 15.01290066666667
$ws.Range("N4").Value = 45.038702
$ws.Range("O4").Value = 0.05872702824751719
$ws.Range("P4").Value = 0.05872702824751717
$ws.Range("Q4").Value = 47.19128172338799
$ws.Range("R4").Value = 424.7215355104919
$ws.Range("S4").Value = 0.0229124175335519
$ws.Range("T4").Value = 0.0229124175335519
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 3.143381999999999
$ws.Range("H5").Value = 9.430145999999999
$ws.Range("I5").Value = 0.3901511487518624
$ws.Range("J5").Value = 0.3901511487518625
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 16.92374966666667
$ws.Range("N5").Value = 50.771249
$ws.Range("O5").Value = 0.06620183179756665
$ws.Range("P5").Value = 0.06620183179756664
$ws.Range("Q5").Value = 53.19781007470599
$ws.Range("R5").Value = 478.7802906723539
$ws.Range("S5").Value = 0.0258287207252982
$ws.Range("T5").Value = 0.0258287207252982
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 3.143381999999999
$ws.Range("H6").Value = 9.430145999999999
$ws.Range("I6").Value = 0.3901511487518624
$ws.Range("J6").Value = 0.3901511487518625
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 12.29188266666667
$ws.Range("N6").Value = 36.875648
$ws.Range("O6").Value = 0.04808302916326276
$ws.Range("P6").Value = 0.04808302916326274
$ws.Range("Q6").Value = 38.63808272051199
$ws.Range("R6").Value = 347.742744484608
$ws.Range("S6").Value = 0.01875964906351626
$ws.Range("T6").Value = 0.01875964906351626
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 4.913449333333333
$ws.Range("H7").Value = 14.740348
$ws.Range("I7").Value = 0.6098488512481375
$ws.Range("J7").Value = 0.6098488512481376
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 86.40747066666665
$ws.Range("N7").Value = 259.222412
$ws.Range("O7").Value = 0.3380062309947018
$ws.Range("P7").Value = 0.3380062309947017
$ws.Range("Q7").Value = 424.5587291421529
$ws.Range("R7").Value = 3821.028562279376
$ws.Range("S7").Value = 0.2061327116868315
$ws.Range("T7").Value = 0.2061327116868314
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 4.913449333333333
$ws.Range("H8").Value = 14.740348
$ws.Range("I8").Value = 0.6098488512481375
$ws.Range("J8").Value = 0.6098488512481376
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 125.002688
$ws.Range("N8").Value = 375.008064
$ws.Range("O8").Value = 0.4889818797969518
$ws.Range("P8").Value = 0.4889818797969517
$ws.Range("Q8").Value = 614.1943740184746
$ws.Range("R8").Value = 5527.749366166272
$ws.Range("S8").Value = 0.2982050376753259
$ws.Range("T8").Value = 0.2982050376753259
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 4.913449333333333
$ws.Range("H9").Value = 14.740348
$ws.Range("I9").Value = 0.6098488512481375
$ws.Range("J9").Value = 0.6098488512481376
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 15.01290066666667
$ws.Range("N9").Value = 45.038702
$ws.Range("O9").Value = 0.05872702824751719
$ws.Range("P9").Value = 0.05872702824751717
$ws.Range("Q9").Value = 73.76512677203289
$ws.Range("R9").Value = 663.8861409482961
$ws.Range("S9").Value = 0.03581461071396528
$ws.Range("T9").Value = 0.03581461071396527
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 4.913449333333333
$ws.Range("H10").Value = 14.740348
$ws.Range("I10").Value = 0.6098488512481375
$ws.Range("J10").Value = 0.6098488512481376
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 16.92374966666667
$ws.Range("N10").Value = 50.771249
$ws.Range("O10").Value = 0.06620183179756665
$ws.Range("P10").Value = 0.06620183179756664
$ws.Range("Q10").Value = 83.15398651718355
$ws.Range("R10").Value = 748.3858786546521
$ws.Range("S10").Value = 0.04037311107226844
$ws.Range("T10").Value = 0.04037311107226844
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 4.913449333333333
$ws.Range("H11").Value = 14.740348
$ws.Range("I11").Value = 0.6098488512481375
$ws.Range("J11").Value = 0.6098488512481376
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 12.29188266666667
$ws.Range("N11").Value = 36.875648
$ws.Range("O11").Value = 0.04808302916326276
$ws.Range("P11").Value = 0.04808302916326274
$ws.Range("Q11").Value = 60.39554269394488
$ws.Range("R11").Value = 543.559884245504
$ws.Range("S11").Value = 0.02932338009974648
$ws.Range("T11").Value = 0.02932338009974648
